$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.092547051789358
$ws.Range("D2").Value = 1.095061633883408
$ws.Range("E2").Value = 1.094370488614895
$ws.Range("F2").Value = 1.105850107209808
$ws.Range("I2").Value = 1.074865683611512
$ws.Range("J2").Value = 1.097367287335862
$ws.Range("K2").Value = 1.097697394051782
$ws.Range("L2").Value = 1.097008004927288
$ws.Range("M2").Value = 1.108458761183434
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.093903339520336
$ws.Range("D3").Value = 1.096186668334008
$ws.Range("E3").Value = 1.095583711910466
$ws.Range("F3").Value = 1.107093062574427
$ws.Range("I3").Value = 1.075402019186792
$ws.Range("J3").Value = 1.098387544671602
$ws.Range("K3").Value = 1.098642669306875
$ws.Range("L3").Value = 1.098041134009828
$ws.Range("M3").Value = 1.109523644434391
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.094779930220307
$ws.Range("D4").Value = 1.096913706579136
$ws.Range("E4").Value = 1.096367522697387
$ws.Range("F4").Value = 1.107896475731043
$ws.Range("I4").Value = 1.075747285865877
$ws.Range("J4").Value = 1.099046155134508
$ws.Range("K4").Value = 1.099252796655645
$ws.Range("L4").Value = 1.098707835193247
$ws.Range("M4").Value = 1.110211258595818
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.095148209904217
$ws.Range("D5").Value = 1.097219133103954
$ws.Range("E5").Value = 1.096696747282983
$ws.Range("F5").Value = 1.108234028007825
$ws.Range("I5").Value = 1.075892012335356
$ws.Range("J5").Value = 1.099322664036074
$ws.Range("K5").Value = 1.099508931349951
$ws.Range("L5").Value = 1.09898768861921
$ws.Range("M5").Value = 1.110499991235379
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.09521003180222
$ws.Range("D6").Value = 1.097270402735816
$ws.Range("E6").Value = 1.096752008680719
$ws.Range("F6").Value = 1.108290692748271
$ws.Range("I6").Value = 1.075916287764306
$ws.Range("J6").Value = 1.099369069448313
$ws.Range("K6").Value = 1.099551916307365
$ws.Range("L6").Value = 1.09903465226605
$ws.Range("M6").Value = 1.110548450860427
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.094784852124898
$ws.Range("D7").Value = 1.096917788565024
$ws.Range("E7").Value = 1.096371922944853
$ws.Range("F7").Value = 1.107900986912526
$ws.Range("I7").Value = 1.075749221370078
$ws.Range("J7").Value = 1.099049851313097
$ws.Range("K7").Value = 1.099256220561724
$ws.Range("L7").Value = 1.098711576284521
$ws.Range("M7").Value = 1.110215117988419
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.0930056292286
$ws.Range("D8").Value = 1.095442039989702
$ws.Range("E8").Value = 1.094780759327748
$ws.Range("F8").Value = 1.10627035056697
$ws.Range("I8").Value = 1.075047310197882
$ws.Range("J8").Value = 1.097712414109123
$ws.Range("K8").Value = 1.09801717275362
$ws.Range("L8").Value = 1.097357531201517
$ws.Range("M8").Value = 1.10881894304293
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.089862402254349
$ws.Range("D9").Value = 1.092834277997748
$ws.Range("E9").Value = 1.091967358149236
$ws.Range("F9").Value = 1.103390189093311
$ws.Range("I9").Value = 1.0737967395714
$ws.Range("J9").Value = 1.095343539088142
$ws.Range("K9").Value = 1.095821958852634
$ws.Range("L9").Value = 1.094957557973353
$ws.Range("M9").Value = 1.106347542973436
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.087761235864152
$ws.Range("D10").Value = 1.091090648293197
$ws.Range("E10").Value = 1.090085071870444
$ws.Range("F10").Value = 1.101465289332866
$ws.Range("I10").Value = 1.072953673272288
$ws.Range("J10").Value = 1.093755902151086
$ws.Range("K10").Value = 1.094350310977555
$ws.Range("L10").Value = 1.093347944010544
$ws.Range("M10").Value = 1.104692213637824
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.086849993315687
$ws.Range("D11").Value = 1.090334375030107
$ws.Range("E11").Value = 1.089268378543
$ws.Range("F11").Value = 1.100630595982354
$ws.Range("I11").Value = 1.072586368100142
$ws.Range("J11").Value = 1.093066397959228
$ws.Range("K11").Value = 1.093711087300125
$ws.Range("L11").Value = 1.092648625495525
$ws.Range("M11").Value = 1.103973553804275
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.08651129754461
$ws.Range("D12").Value = 1.090053266462169
$ws.Range("E12").Value = 1.088964769250033
$ws.Range("F12").Value = 1.100320368829791
$ws.Range("I12").Value = 1.072449593642548
$ws.Range("J12").Value = 1.092809973427753
$ws.Range("K12").Value = 1.093473348106866
$ws.Range("L12").Value = 1.092388510457846
$ws.Range("M12").Value = 1.103706322940539
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.086583958973783
$ws.Range("D13").Value = 1.090113574096144
$ws.Range("E13").Value = 1.089029906004497
$ws.Range("F13").Value = 1.100386922049509
$ws.Range("I13").Value = 1.072478947699801
$ws.Range("J13").Value = 1.092864991549394
$ws.Range("K13").Value = 1.09352435774962
$ws.Range("L13").Value = 1.092444322247127
$ws.Range("M13").Value = 1.103763657999669
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.086822001133554
$ws.Range("D14").Value = 1.090311142510542
$ws.Range("E14").Value = 1.089243287294791
$ws.Range("F14").Value = 1.100604956308065
$ws.Range("I14").Value = 1.072575069251658
$ws.Range("J14").Value = 1.09304520821746
$ws.Range("K14").Value = 1.093691441909856
$ws.Range("L14").Value = 1.092627131615336
$ws.Range("M14").Value = 1.103951470330992
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.086968637475198
$ws.Range("D15").Value = 1.090432844965479
$ws.Range("E15").Value = 1.089374724871778
$ws.Range("F15").Value = 1.100739269770677
$ws.Range("I15").Value = 1.072634247652136
$ws.Range("J15").Value = 1.093156204181447
$ws.Range("K15").Value = 1.09379434769675
$ws.Range("L15").Value = 1.092739719053432
$ws.Range("M15").Value = 1.10406714933047
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.087821681539685
$ws.Range("D16").Value = 1.091140812542665
$ws.Range("E16").Value = 1.090139237907869
$ws.Range("F16").Value = 1.101520659515507
$ws.Range("I16").Value = 1.072978002436917
$ws.Range("J16").Value = 1.09380161876813
$ws.Range("K16").Value = 1.094392691869209
$ws.Range("L16").Value = 1.093394305670273
$ws.Range("M16").Value = 1.104739868496713
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.088356388289691
$ws.Range("D17").Value = 1.091584558910454
$ws.Range("E17").Value = 1.090618351215837
$ws.Range("F17").Value = 1.102010480304639
$ws.Range("I17").Value = 1.073193025974294
$ws.Range("J17").Value = 1.094205919135191
$ws.Range("K17").Value = 1.094767481634902
$ws.Range("L17").Value = 1.093804279423579
$ws.Range("M17").Value = 1.105161337961808
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.088668136726491
$ws.Range("D18").Value = 1.09184326639853
$ws.Range("E18").Value = 1.090897650949365
$ws.Range("F18").Value = 1.102296069114782
$ws.Range("I18").Value = 1.073318228469438
$ws.Range("J18").Value = 1.094441543503309
$ws.Range("K18").Value = 1.094985898389372
$ws.Range("L18").Value = 1.094043184233238
$ws.Range("M18").Value = 1.105406991789388
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.088774411831224
$ws.Range("D19").Value = 1.091931458354769
$ws.Range("E19").Value = 1.090992858129789
$ws.Range("F19").Value = 1.102393428092579
$ws.Range("I19").Value = 1.073360882540902
$ws.Range("J19").Value = 1.094521852009796
$ws.Range("K19").Value = 1.09506034045078
$ws.Range("L19").Value = 1.094124606416493
$ws.Range("M19").Value = 1.105490722629401
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.0882990335264
$ws.Range("D20").Value = 1.091536961802311
$ws.Range("E20").Value = 1.090566963370952
$ws.Range("F20").Value = 1.10195793913837
$ws.Range("I20").Value = 1.073169978458848
$ws.Range("J20").Value = 1.094162561962379
$ws.Range("K20").Value = 1.094727290106424
$ws.Range("L20").Value = 1.093760316525503
$ws.Range("M20").Value = 1.105116137154266
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.086751909762752
$ws.Range("D21").Value = 1.090252968953618
$ws.Range("E21").Value = 1.089180458848225
$ws.Range("F21").Value = 1.100540755794961
$ws.Range("I21").Value = 1.072546773283906
$ws.Range("J21").Value = 1.092992147548364
$ws.Range("K21").Value = 1.093642248187398
$ws.Range("L21").Value = 1.092573308705631
$ws.Range("M21").Value = 1.103896172283291
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.085777897531628
$ws.Range("D22").Value = 1.089444541357288
$ws.Range("E22").Value = 1.088307241991662
$ws.Range("F22").Value = 1.099648644371218
$ws.Range("I22").Value = 1.072152965434763
$ws.Range("J22").Value = 1.092254454344442
$ws.Range("K22").Value = 1.092958283861301
$ws.Range("L22").Value = 1.091824921922994
$ws.Range("M22").Value = 1.103127459954871
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.086294362569606
$ws.Range("D23").Value = 1.089873212668504
$ws.Range("E23").Value = 1.088770291443004
$ws.Range("F23").Value = 1.100121672717036
$ws.Range("I23").Value = 1.072361918425668
$ws.Range("J23").Value = 1.092645692262647
$ws.Range("K23").Value = 1.093321034030585
$ws.Range("L23").Value = 1.092221853495737
$ws.Range("M23").Value = 1.103535128883843
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.088324950097977
$ws.Range("D24").Value = 1.091558469264249
$ws.Range("E24").Value = 1.090590183816961
$ws.Range("F24").Value = 1.101981680586032
$ws.Range("I24").Value = 1.073180393309649
$ws.Range("J24").Value = 1.094182153811536
$ws.Range("K24").Value = 1.094745451520452
$ws.Range("L24").Value = 1.093780182164645
$ws.Range("M24").Value = 1.105136562016215
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.090675981354333
$ws.Range("D25").Value = 1.093509334200746
$ws.Range("E25").Value = 1.092695849543978
$ws.Range("F25").Value = 1.104135606906983
$ws.Range("I25").Value = 1.0741216807572
$ws.Range("J25").Value = 1.095957411762091
$ws.Range("K25").Value = 1.096390899718364
$ws.Range("L25").Value = 1.095579688785507
$ws.Range("M25").Value = 1.106987804904661
